# Insert a new data row at row 366 (pushing the existing rows 366-508 down
# to 367-509), then populate the newly inserted row with the new record's
# values. This reproduces the target diff, which is a single new row
# inserted in the middle of the "Apio" price table (dimension grows from
# A1:R508 to A1:R509).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above current row 366; Excel shifts rows 366..508
# down to 367..509 and extends the used range automatically.
$ws.Rows.Item(366).Insert()

# Populate the newly inserted row 366 with the new record.
$ws.Cells.Item(366, 1).Value  = 4
$ws.Cells.Item(366, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(366, 3).Value  = "Los Lagos"
$ws.Cells.Item(366, 4).Value  = 45141
$ws.Cells.Item(366, 5).Value  = 10
$ws.Cells.Item(366, 6).Value  = 100112017
$ws.Cells.Item(366, 7).Value  = "Apio"
$ws.Cells.Item(366, 8).Value  = "Americana (o)"
$ws.Cells.Item(366, 9).Value  = "Primera"
$ws.Cells.Item(366, 10).Value = 25
$ws.Cells.Item(366, 11).Value = 11000
$ws.Cells.Item(366, 12).Value = 11000
$ws.Cells.Item(366, 13).Value = 11000
$ws.Cells.Item(366, 14).Value = "`$/docena de matas"
$ws.Cells.Item(366, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(366, 16).Value = 1833
$ws.Cells.Item(366, 17).Value = 6
$ws.Cells.Item(366, 18).Value = "Hortaliza"

# Apply the same date number format style used by the rest of column D so
# the new date cell renders/serialises identically to its neighbours.
$ws.Cells.Item(366, 4).NumberFormat = $ws.Cells.Item(367, 4).NumberFormat
